$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "44 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "4|    |"
$t.Cell(1,2).Range.Text = "49 x 10" + [char]11 + "  1    0" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "9|    |"
$t.Cell(1,3).Range.Text = "60 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "0|    |"
$t.Cell(2,1).Range.Text = "58 x 69" + [char]11 + "  6    9" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "8|    |"
$t.Cell(2,2).Range.Text = "93 x 84" + [char]11 + "  8    4" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "3|    |"
$t.Cell(2,3).Range.Text = "80 x 29" + [char]11 + "  2    9" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "0|    |"
$t.Cell(3,1).Range.Text = "38 x 49" + [char]11 + "  4    9" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "8|    |"
$t.Cell(3,2).Range.Text = "26 x 96" + [char]11 + "  9    6" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "6|    |"
$t.Cell(3,3).Range.Text = "70 x 64" + [char]11 + "  6    4" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "0|    |"
$t.Cell(4,1).Range.Text = "36 x 56" + [char]11 + "  5    6" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "6|    |"
$t.Cell(4,2).Range.Text = "18 x 12" + [char]11 + "  1    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "8|    |"
$t.Cell(4,3).Range.Text = "41 x 58" + [char]11 + "  5    8" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "1|    |"
$t.Cell(5,1).Range.Text = "74 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "4|    |"
$t.Cell(5,2).Range.Text = "17 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"
$t.Cell(5,3).Range.Text = "11 x 11" + [char]11 + "  1    1" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11
